$d = $word.ActiveDocument

foreach ($p in $d.Paragraphs) {
    $rng = $p.Range
    if ($rng.Font.HighlightColorIndex) {
        $rng.Font.HighlightColorIndex = 0
    }
}
